# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff", and the handoff/generation timestamps were
# refreshed. Update every sheet that carries this status/timestamp
# pair, then widen the "Status" column (and its mirrored columns on
# the Overview tab) so the longer text isn't truncated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps -------------------------------------------------
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" both advance from 06:57:05 to 06:57:35.
$wsOverview.Range("G2").Value = "2016-08-25 06:57:35"
$wsDeDe.Range("H2").Value     = "2016-08-25 06:57:35"

# zh-cn's "Latest Handoff Datetime" advances from 06:56:57 to 06:57:30.
$wsZhCn.Range("H2").Value = "2016-08-25 06:57:30"

# --- Column widths: widen the Status column(s) for "Ready for handoff" ---
# (13.4101845877511 -> 17.2159881591797 character-width units once stored;
# ColumnWidth is expressed relative to the workbook's standard font.)
$newColumnWidth = 16.3826548258464

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth     = $newColumnWidth  # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = $newColumnWidth  # column C (Status)
